# Generate Report for Handoff
# - Flip status from "Handed back: in sync with en-US" to "Ready for handoff"
#   on the Overview sheet (zh-cn/de-de status columns) and on each
#   language sheet's "Status" column.
# - Bump the associated timestamps to reflect the new handoff generation time.
# - Re-fit the "Status" column(s) now that the text is shorter.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-31 13:15:17"

# Narrow the zh-cn / de-de status columns (E, F) to fit the shorter text.
$overview.Columns.Item(5).ColumnWidth = 16.33
$overview.Columns.Item(6).ColumnWidth = 16.33

# ---- zh-cn sheet ----
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-31 13:14:58"
$zhcn.Columns.Item(3).ColumnWidth = 16.33

# ---- de-de sheet ----
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-31 13:15:17"
$dede.Columns.Item(3).ColumnWidth = 16.33
